# Regenerate merged AHB files
#
# 1. Rename the "old"/"new" header-column convention to the concrete
#    file-version labels (FV2210 / FV2304).
# 2. Turn A1:U66 into a genuine Excel Table ("Table1") with banded rows.
# 3. Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header cells (row 1).
#    Columns A-J: "<Name>_old" -> "<Name>_FV2210"
#    Column  K:   "diff"       -> unchanged
#    Columns L-U: "<Name>_new" -> "<Name>_FV2304"
# ---------------------------------------------------------------------------
$headerRenames = @{
    "A1" = "Segmentname_FV2210"
    "B1" = "Segmentgruppe_FV2210"
    "C1" = "Segment_FV2210"
    "D1" = "Datenelement_FV2210"
    "E1" = "Segment ID_FV2210"
    "F1" = "Code_FV2210"
    "G1" = "Qualifier_FV2210"
    "H1" = "Beschreibung_FV2210"
    "I1" = "Bedingungsausdruck_FV2210"
    "J1" = "Bedingung_FV2210"
    "L1" = "Segmentname_FV2304"
    "M1" = "Segmentgruppe_FV2304"
    "N1" = "Segment_FV2304"
    "O1" = "Datenelement_FV2304"
    "P1" = "Segment ID_FV2304"
    "Q1" = "Code_FV2304"
    "R1" = "Qualifier_FV2304"
    "S1" = "Beschreibung_FV2304"
    "T1" = "Bedingungsausdruck_FV2304"
    "U1" = "Bedingung_FV2304"
}

foreach ($addr in $headerRenames.Keys) {
    $ws.Range($addr).Value = $headerRenames[$addr]
}

# ---------------------------------------------------------------------------
# 2. Turn A1:U66 into a real table ("Table1") with banded rows.
#    The header row already carries bold / fill / border formatting. If that
#    formatting is left in place while the ListObject is created, the engine
#    captures it into a brand-new header-row dxf (headerRowDxfId) - something
#    the source workbook's table does NOT have. So: stash the header's
#    current look, strip it down to the default style, create the table,
#    then paste the original look back from the stash (format-only paste, so
#    no extra style entries get synthesized).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$styleStash = $ws.Range("A68:U68")

$headerRange.Copy()
$styleStash.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.Style = "Normal"

$tableRange = $ws.Range("A1:U66")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

$styleStash.Copy()
$headerRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$styleStash.Clear() | Out-Null

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
